$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the text in J6 (Reason of the Change) and K6 (Changes Done(In Details)) ---
$ws.Range("J6").Value = "1.DB status check for hana DB had to be written.`n2.DB status check for db2 DB had to be written.`n3.When changing the commands to execute using ""sudo bash"", an extra ""/"" after ""'"""
$ws.Range("K6").Value = "1.For Hana DB, the ""sapcontrol"" command is run with the function ""GetProcessList"" to check the status of the hana db and output is printed based on the output of the command.`n2.For DB2 DB, we check if ""db2sysc"" process is running, if it is running it means that the DB is up. The output is printed based on the output of the command.`n3.The ""/"" is removed."

# --- 2. Add a new column L ("Applicable") ---
$ws.Range("L5").Value = "Applicable"
$ws.Range("L6").Value = "Generic/Enhancement(1,2)/Error(3)"

# Match formatting of the neighbouring column K for the new column L
$ws.Range("K5").Copy() | Out-Null
$ws.Range("L5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("K6").Copy() | Out-Null
$ws.Range("L6").PasteSpecial(-4122) | Out-Null

$ws.Range("K7:K16").Copy() | Out-Null
$ws.Range("L7:L16").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = $false

# Column width for the new column (closest achievable to 23.1796875 given rounding)
$ws.Columns.Item(12).ColumnWidth = 22.3

# --- 3. Update the view: scroll position + selection ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 3
$ws.Range("L6:L16").Select() | Out-Null
